# Turnos assimétricos (de dia para dia a existencia de turnos pode variar)
#
# Remove the three "Coordinator" rows (Nome34, Nome40, Nome44 — that team
# no longer exists) and append a batch of newly-collected volunteers
# (Nome52 .. Nome94, with gaps — not every numbered name was filled in),
# then re-sort the whole table by Name and refresh the AutoFilter so it
# covers the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the rows that belonged to people no longer in the roster ---
# (delete from the bottom up so earlier row numbers stay valid)
$ws.Rows.Item(45).Delete()   # Nome44 / Coordinator
$ws.Rows.Item(41).Delete()   # Nome40 / Coordinator
$ws.Rows.Item(35).Delete()   # Nome34 / Coordinator

# --- 2. Append the newly added people at the bottom of the sheet ---
$availability = "[1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1]"

$newPeople = @(
    @("Nome52", "WebDev"),
    @("Nome53", "Speakers"),
    @("Nome54", "WebDev"),
    @("Nome55", "Speakers"),
    @("Nome56", "WebDev"),
    @("Nome57", "WebDev"),
    @("Nome58", "Marketing"),
    @("Nome59", "Marketing"),
    @("Nome61", "Marketing"),
    @("Nome63", "Speakers"),
    @("Nome65", "Marketing"),
    @("Nome66", "Marketing"),
    @("Nome67", "Speakers"),
    @("Nome69", "Marketing"),
    @("Nome74", "Business"),
    @("Nome75", "Volunteer"),
    @("Nome76", "Volunteer"),
    @("Nome77", "Volunteer"),
    @("Nome78", "Business"),
    @("Nome83", "Business"),
    @("Nome86", "Logistics"),
    @("Nome88", "Logistics"),
    @("Nome89", "Logistics"),
    @("Nome92", "Logistics"),
    @("Nome94", "Business")
)

$r = 50
foreach ($person in $newPeople) {
    $ws.Cells.Item($r, 1).Value = $person[0]
    $ws.Cells.Item($r, 2).Value = $person[1]
    $ws.Cells.Item($r, 3).Value = $availability
    $r = $r + 1
}

# --- 3. Re-sort the whole table (header + 73 people) by Name ---
$fullRange = $ws.Range("A1:C74")
$fullRange.Sort($ws.Range("A1"))

# --- 4. Refresh the AutoFilter over the new, bigger range ---
$fullRange.AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Folha1!`$A`$1:`$C`$74")
$filterName.Visible = $false

# --- 5. Leave the selection where the user left it after scrolling down ---
$ws.Range("C79").Select()
